# auto increment public id and fix unit tests
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Parent Public ID) - renumber statement references
$ws.Range("A2").Value = "statement-01"
$ws.Range("A3").Value = "statement-01"
$ws.Range("A4").Value = "statement-01"
$ws.Range("A5").Value = "statement-02"
$ws.Range("A6").Value = "statement-02"

# Column B (Public ID) - clear auto-increment id values
$ws.Range("B2").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("B4").Value = ""
$ws.Range("B5").Value = ""
$ws.Range("B6").Value = ""
